$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '62.915.64'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -2.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.125.45'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.67'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -2.44%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '136.28'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -4.58%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.125.02'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.514'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.56%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.59%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '5.23'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -2.90%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.454'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  -3.03%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000242'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -5.08%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '34.01'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -3.25%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.638.59'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.50%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '62.972.47'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.125.20'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.55%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.62'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -3.93%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '470.33'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.09'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -3.53%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.696'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -1.96%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.66'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.28'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.41%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.92'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -3.70%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -1.89%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.87'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -5.60%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.09'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +2.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.84'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -4.29%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '26.56'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.92%  '
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -4.30%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.53'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.19%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -3.32%  '
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.74'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.84%  '
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '51.96'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0₃0677'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -11.92%  '
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -2.16%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '416.02'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -6.74%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.915.06'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.17'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.67'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -10.92%  '
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -6.15%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.259'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +0.10%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.09'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -5.59%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '25.30'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -2.62%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.25'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -7.82%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '120.28'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.38%  '
